# Atualização de bases das ligas, do dia: 15-04-2024 às 22:35
#
# - Swap match records for existing rows 83/84 (reorders the two October
#   2023 Cavalry FC vs Forge FC / Atletico Ottawa vs Pacific FC CA fixtures).
# - Fill in the previously partial row 89 with the full set of odds columns.
# - Append two brand-new fixture rows (90, 91) for the April 2024 round.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row83Data = @{
    2 = 7301364   # B83
    3 = "Canada Premier League"   # C83
    4 = "Canada Premier League"   # D83
    5 = 45206.75   # E83
    6 = "Forge FC"   # F83
    7 = "Atletico Ottawa"   # G83
    8 = 0   # H83
    9 = 1   # I83
    10 = "A"   # J83
    11 = 1.8   # K83
    12 = 3.6   # L83
    13 = 3.5   # M83
    14 = 1.533   # N83
    15 = 3.8   # O83
    16 = 5   # P83
    17 = -1   # Q83
    18 = 1.975   # R83
    19 = 1.825   # S83
    20 = 2.5   # T83
    21 = 1.9   # U83
    22 = 1.9   # V83
    23 = -1   # W83
    24 = -1   # X83
    25 = 4   # Y83
    26 = -1   # Z83
    27 = 0.825   # AA83
    28 = -1   # AB83
    29 = 0.8999999999999999   # AC83
}
$row84Data = @{
    2 = 6227884   # B84
    3 = "Canada Premier League"   # C84
    4 = "Canada Premier League"   # D84
    5 = 45206.75   # E84
    6 = "Cavalry FC"   # F84
    7 = "Pacific FC CA"   # G84
    8 = 3   # H84
    9 = 0   # I84
    10 = "H"   # J84
    11 = 2.25   # K84
    12 = 3.1   # L84
    13 = 2.875   # M84
    14 = 2.05   # N84
    15 = 3.2   # O84
    16 = 3.2   # P84
    17 = -0.25   # Q84
    18 = 1.825   # R84
    19 = 1.975   # S84
    20 = 2.5   # T84
    21 = 1.825   # U84
    22 = 1.975   # V84
    23 = 1.05   # W84
    24 = -1   # X84
    25 = -1   # Y84
    26 = 0.825   # Z84
    27 = -1   # AA84
    28 = 0.825   # AB84
    29 = -1   # AC84
}
$row89Data = @{
    1 = 87   # A89
    2 = 7802934   # B89
    3 = "Canada Premier League"   # C89
    4 = "Canada Premier League"   # D89
    5 = 45395.58333333334   # E89
    6 = "Atletico Ottawa"   # F89
    7 = "York United FC"   # G89
    8 = 2   # H89
    9 = 1   # I89
    10 = "H"   # J89
    11 = 2.875   # K89
    12 = 3.4   # L89
    13 = 2.1   # M89
    14 = 2.4   # N89
    15 = 3.4   # O89
    16 = 2.5   # P89
    17 = 0   # Q89
    18 = 1.85   # R89
    19 = 1.95   # S89
    20 = 2.5   # T89
    21 = 2   # U89
    22 = 1.8   # V89
    23 = 1.4   # W89
    24 = -1   # X89
    25 = -1   # Y89
    26 = 0.8500000000000001   # Z89
    27 = -1   # AA89
    28 = 1   # AB89
    29 = -1   # AC89
}
$row90Data = @{
    1 = 88   # A90
    2 = 7802874   # B90
    3 = "Canada Premier League"   # C90
    4 = "Canada Premier League"   # D90
    5 = 45395.70833333334   # E90
    6 = "Forge FC"   # F90
    7 = "Cavalry FC"   # G90
    8 = 2   # H90
    9 = 1   # I90
    10 = "H"   # J90
    11 = 2.2   # K90
    12 = 3.6   # L90
    13 = 2.6   # M90
    14 = 2.25   # N90
    15 = 3.6   # O90
    16 = 2.55   # P90
    17 = 0   # Q90
    18 = 1.775   # R90
    19 = 2.025   # S90
    20 = 2.5   # T90
    21 = 1.95   # U90
    22 = 1.85   # V90
    23 = 1.25   # W90
    24 = -1   # X90
    25 = -1   # Y90
    26 = 0.7749999999999999   # Z90
    27 = -1   # AA90
    28 = 0.95   # AB90
    29 = -1   # AC90
}
$row91Data = @{
    1 = 89   # A91
    2 = 7803361   # B91
    3 = "Canada Premier League"   # C91
    4 = "Canada Premier League"   # D91
    5 = 45395.83333333334   # E91
    6 = "Pacific FC CA"   # F91
    7 = "HFX Wanderers"   # G91
    8 = 1   # H91
    9 = 0   # I91
    10 = "H"   # J91
    11 = 1.833   # K91
    12 = 3.5   # L91
    13 = 3.5   # M91
    14 = 1.95   # N91
    15 = 3.25   # O91
    16 = 3.4   # P91
    17 = -0.25   # Q91
    18 = 1.775   # R91
    19 = 2.025   # S91
    20 = 2.25   # T91
    21 = 1.875   # U91
    22 = 1.925   # V91
    23 = 0.95   # W91
    24 = -1   # X91
    25 = -1   # Y91
    26 = 0.7749999999999999   # Z91
    27 = -1   # AA91
    28 = -1   # AB91
    29 = 0.925   # AC91
}

function Set-RowValues($Worksheet, $RowNumber, $Values) {
    foreach ($colNum in $Values.Keys) {
        $Worksheet.Cells.Item($RowNumber, [int]$colNum).Value = $Values[$colNum]
    }
}

# --- Rows 83 & 84: the two fixtures' data (columns B:AC) are swapped;
#     column A (the sequential id 81 / 82) stays attached to its row. ---
Set-RowValues $ws 83 $row83Data
Set-RowValues $ws 84 $row84Data

# --- Row 89: refresh with the corrected/complete record (adds the
#     previously-missing FTHG/FTAG/FTR and PL_AhOver/PL_AhUnder cells). ---
Set-RowValues $ws 89 $row89Data

# --- Rows 90 & 91: brand-new fixture rows. Clone row 89's formatting
#     first (bold/boxed id in col A, date format in col E) then fill values. ---
$ws.Range("A89:AC89").Copy()
$ws.Range("A90:AC90").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A89:AC89").Copy()
$ws.Range("A91:AC91").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Set-RowValues $ws 90 $row90Data
Set-RowValues $ws 91 $row91Data
